$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'65.043.71"
$ws.Cells.Item(2, 5).Value = "  -0.03%  "
$ws.Cells.Item(3, 4).Value = "'3.519.24"
$ws.Cells.Item(3, 5).Value = "  -1.54%  "
$ws.Cells.Item(4, 5).Value = "  +0.09%  "
$ws.Cells.Item(5, 4).Value = "'592.72"
$ws.Cells.Item(5, 5).Value = "  -1.61%  "
$ws.Cells.Item(6, 4).Value = "'134.14"
$ws.Cells.Item(6, 5).Value = "  -1.27%  "
$ws.Cells.Item(7, 4).Value = "'3.517.78"
$ws.Cells.Item(7, 5).Value = "  -1.56%  "
$ws.Cells.Item(8, 5).Value = "  +0.04%  "
$ws.Cells.Item(9, 4).Value = "'0.490"
$ws.Cells.Item(9, 5).Value = "  -0.72%  "
$ws.Cells.Item(10, 4).Value = "'0.124"
$ws.Cells.Item(10, 5).Value = "  +0.99%  "
$ws.Cells.Item(11, 4).Value = "'7.13"
$ws.Cells.Item(11, 5).Value = "  +2.51%  "
$ws.Cells.Item(12, 5).Value = "  -0.31%  "
$ws.Cells.Item(13, 4).Value = "'4.122.04"
$ws.Cells.Item(13, 5).Value = "  -1.40%  "
$ws.Cells.Item(14, 5).Value = "  +1.57%  "
$ws.Cells.Item(15, 4).Value = "'0.0000181"
$ws.Cells.Item(15, 5).Value = "  -0.83%  "
$ws.Cells.Item(17, 4).Value = "'3.525.35"
$ws.Cells.Item(17, 5).Value = "  -1.35%  "
$ws.Cells.Item(18, 4).Value = "'65.047.69"
$ws.Cells.Item(18, 5).Value = "  -0.15%  "
$ws.Cells.Item(19, 4).Value = "'10.09"
$ws.Cells.Item(19, 5).Value = "  -0.30%  "
$ws.Cells.Item(20, 4).Value = "'14.35"
$ws.Cells.Item(20, 5).Value = "  -0.49%  "
$ws.Cells.Item(21, 4).Value = "'5.67"
$ws.Cells.Item(21, 5).Value = "  -3.31%  "
$ws.Cells.Item(22, 4).Value = "'392.39"
$ws.Cells.Item(22, 5).Value = "  +0.96%  "
$ws.Cells.Item(23, 4).Value = "'0.578"
$ws.Cells.Item(23, 5).Value = "  -0.38%  "
$ws.Cells.Item(24, 4).Value = "'3.664.13"
$ws.Cells.Item(24, 5).Value = "  -1.45%  "
$ws.Cells.Item(25, 4).Value = "'74.57"
$ws.Cells.Item(25, 5).Value = "  +0.41%  "
$ws.Cells.Item(26, 5).Value = "  -0.07%  "
$ws.Cells.Item(27, 5).Value = "  -5.04%  "
$ws.Cells.Item(28, 4).Value = "'1.59"
$ws.Cells.Item(28, 5).Value = "  +8.74%  "
$ws.Cells.Item(29, 4).Value = "'7.66"
$ws.Cells.Item(29, 5).Value = "  -1.00%  "
$ws.Cells.Item(30, 4).Value = "'0.999"
$ws.Cells.Item(30, 5).Value = "  +0.29%  "
$ws.Cells.Item(31, 5).Value = "  -1.66%  "
$ws.Cells.Item(32, 4).Value = "'8.31"
$ws.Cells.Item(32, 5).Value = "  -1.10%  "
$ws.Cells.Item(33, 4).Value = "'3.528.55"
$ws.Cells.Item(33, 5).Value = "  -1.49%  "
$ws.Cells.Item(34, 4).Value = "'24.09"
$ws.Cells.Item(34, 5).Value = "  +0.26%  "
$ws.Cells.Item(35, 5).Value = "  -0.02%  "
$ws.Cells.Item(36, 5).Value = "  -0.25%  "
$ws.Cells.Item(37, 4).Value = "'5.27"
$ws.Cells.Item(37, 5).Value = "  +4.68%  "
$ws.Cells.Item(38, 5).Value = "  +1.15%  "
$ws.Cells.Item(39, 4).Value = "'6.94"
$ws.Cells.Item(39, 5).Value = "  -0.07%  "
$ws.Cells.Item(40, 4).Value = "'168.04"
$ws.Cells.Item(40, 5).Value = "  -0.72%  "
$ws.Cells.Item(41, 4).Value = "'0.0805"
$ws.Cells.Item(41, 5).Value = "  -0.45%  "
$ws.Cells.Item(42, 4).Value = "'0.821"
$ws.Cells.Item(42, 5).Value = "  -0.80%  "
$ws.Cells.Item(43, 5).Value = "  +4.33%  "
$ws.Cells.Item(44, 5).Value = "  +0.65%  "
$ws.Cells.Item(45, 4).Value = "'25.65"
$ws.Cells.Item(45, 5).Value = "  -5.50%  "
$ws.Cells.Item(46, 5).Value = "  +0.14%  "
$ws.Cells.Item(47, 4).Value = "'4.43"
$ws.Cells.Item(47, 5).Value = "  -1.32%  "
$ws.Cells.Item(48, 4).Value = "'1.66"
$ws.Cells.Item(48, 5).Value = "  +0.13%  "
$ws.Cells.Item(49, 4).Value = "'6.90"
$ws.Cells.Item(49, 5).Value = "  -0.58%  "
$ws.Cells.Item(50, 4).Value = "'2.418.21"
$ws.Cells.Item(50, 5).Value = "  -3.61%  "
$ws.Cells.Item(51, 4).Value = "'0.904"
